$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.832.17"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.444.94"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "484.10"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.61"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.55%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.451.34"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  +8.70%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.878.84"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "55.902.20"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.81"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.457.56"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.48"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.99"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.88"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.78"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.25"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.567.75"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.63"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0772"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.42%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.69"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.09"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.11"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.67"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.844"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.73"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("E40").Value = "  +7.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0545"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.31"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "259.36"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0919"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.93%  "
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +11.68%  "
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.36"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.858.19"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.02%  "
